$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '71.170.20'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '3.842.99'
$ws.Range("E3").Value = '  +0.64%  '

Set-TextCell "D4" '0.999'
$ws.Range("E4").Value = '  -0.09%  '

Set-TextCell "D5" '707.29'
$ws.Range("E5").Value = '  +0.94%  '

Set-TextCell "D6" '172.96'
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").Value = '3.842.61'
$ws.Range("E7").Value = '  +0.64%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("E10").Value = '  +0.55%  '

Set-TextCell "D11" '7.29'
$ws.Range("E11").Value = '  +0.52%  '

$ws.Range("E12").Value = '  -0.01%  '

$ws.Range("E13").Value = '  +0.44%  '

Set-TextCell "D14" '36.66'
$ws.Range("E14").Value = '  +1.10%  '

$ws.Range("D15").Value = '4.492.00'
$ws.Range("E15").Value = '  +0.62%  '

$ws.Range("D16").Value = '3.848.94'
$ws.Range("E16").Value = '  +0.45%  '

$ws.Range("D17").Value = '71.109.17'
$ws.Range("E17").Value = '  +0.20%  '

$ws.Range("E18").Value = '  +0.05%  '

Set-TextCell "D20" '17.37'
$ws.Range("E20").Value = '  -2.27%  '

Set-TextCell "D21" '10.75'
$ws.Range("E21").Value = '  -3.85%  '

Set-TextCell "D22" '494.84'
$ws.Range("E22").Value = '  +3.19%  '

Set-TextCell "D23" '0.725'
$ws.Range("E23").Value = '  +1.79%  '

Set-TextCell "D24" '85.02'
$ws.Range("E24").Value = '  +1.38%  '

$ws.Range("E25").Value = '  +2.55%  '

Set-TextCell "D26" '10.65'
$ws.Range("E26").Value = '  +1.74%  '

$ws.Range("E27").Value = '  -1.39%  '

Set-TextCell "D28" '2.11'
$ws.Range("E28").Value = '  -2.64%  '

$ws.Range("E29").Value = '  +2.43%  '

$ws.Range("E30").Value = '  -0.08%  '

Set-TextCell "D31" '7.51'
$ws.Range("E31").Value = '  -0.13%  '

$ws.Range("E32").Value = '  -1.74%  '

Set-TextCell "D33" '29.50'
$ws.Range("E33").Value = '  -0.22%  '

Set-TextCell "D34" '0.179'
$ws.Range("E34").Value = '  -3.13%  '

Set-TextCell "D35" '9.22'
$ws.Range("E35").Value = '  -0.32%  '

$ws.Range("D36").Value = '3.817.16'
$ws.Range("E36").Value = '  +1.26%  '

Set-TextCell "D37" '0.998'
$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("E39").Value = '  +5.51%  '

$ws.Range("E40").Value = '  +5.31%  '

Set-TextCell "D41" '6.03'
$ws.Range("E41").Value = '  -0.02%  '

Set-TextCell "D42" '3.37'
$ws.Range("E42").Value = '  -1.32%  '

$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("E45").Value = '  -2.37%  '

Set-TextCell "D46" '163.50'
$ws.Range("E46").Value = '  +0.45%  '

Set-TextCell "D47" '48.67'
$ws.Range("E47").Value = '  -0.54%  '

$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell "D48" '1.39'
$ws.Range("E48").Value = '  +0.48%  '

$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell "D49" '415.81'
$ws.Range("E49").Value = '  +1.32%  '

Set-TextCell "D50" '8.61'
$ws.Range("E50").Value = '  +0.64%  '

Set-TextCell "D51" '0.298'
$ws.Range("E51").Value = '  -0.98%  '
